$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4262
$ws.Range("I74").Value = 3649.5
$ws.Range("J74").Value = 4874.5
$ws.Range("K74").Value = 3649.5
$ws.Range("L74").Value = 4874.5
$ws.Range("M74").Value = -2713.5
$ws.Range("N74").Value = -6746.5

$ws.Range("H76").Value = 4375.25
$ws.Range("I76").Value = 4003
$ws.Range("J76").Value = 4499.3335
$ws.Range("K76").Value = 4003
$ws.Range("L76").Value = 4499.3335
$ws.Range("M76").Value = -3688
$ws.Range("N76").Value = -5129.3335

$ws.Range("H77").Value = 4262
$ws.Range("I77").Value = 3649.5
$ws.Range("J77").Value = 4874.5
$ws.Range("K77").Value = 18247.5
$ws.Range("L77").Value = 24372.5
$ws.Range("M77").Value = -13567.5
$ws.Range("N77").Value = -33732.5

$ws.Range("H79").Value = 4375.25
$ws.Range("I79").Value = 4003
$ws.Range("J79").Value = 4499.3335
$ws.Range("K79").Value = 4003
$ws.Range("L79").Value = 4499.3335
$ws.Range("M79").Value = -2911
$ws.Range("N79").Value = -6683.3335

$ws.Range("H116").Value = 2648.3333
$ws.Range("I116").Value = 2422.2222
$ws.Range("J116").Value = 3326.6667
$ws.Range("K116").Value = 2422.2222
$ws.Range("L116").Value = 3326.6667
$ws.Range("M116").Value = 1019.7778
$ws.Range("N116").Value = -10210.6667

$ws.Range("H135").Value = 996.89795
$ws.Range("I135").Value = 949.1163
$ws.Range("J135").Value = 1339.3334
$ws.Range("K135").Value = 8542.046700000001
$ws.Range("L135").Value = 12054.0006
$ws.Range("M135").Value = -6007.046700000001
$ws.Range("N135").Value = -17124.0006

$ws.Range("H137").Value = 10868.429
$ws.Range("I137").Value = 959.5294
$ws.Range("J137").Value = 26182.182
$ws.Range("K137").Value = 2878.5882
$ws.Range("L137").Value = 78546.546
$ws.Range("M137").Value = -328.5882000000001
$ws.Range("N137").Value = -83646.546

$ws.Range("H138").Value = 3971976.5
$ws.Range("I138").Value = 11907188
$ws.Range("J138").Value = 4370.768
$ws.Range("K138").Value = 35721564
$ws.Range("L138").Value = 13112.304
$ws.Range("M138").Value = -35716424
$ws.Range("N138").Value = -23392.304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20916.453
$ws.Range("I32").Value = 20991.607
$ws.Range("J32").Value = 19000
$ws.Range("K32").Value = 20991.607
$ws.Range("L32").Value = 19000
$ws.Range("M32").Value = -20704.607
$ws.Range("N32").Value = -19574

$ws.Range("H63").Value = 62503270
$ws.Range("I63").Value = 83336370
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 83336370
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -83335684
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 62503270
$ws.Range("I66").Value = 83336370
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 416681850
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -416678418
$ws.Range("N66").Value = -26864

$ws.Range("H122").Value = 13578.167
$ws.Range("I122").Value = 16607.334
$ws.Range("J122").Value = 4490.6665
$ws.Range("K122").Value = 49822.00199999999
$ws.Range("L122").Value = 13471.9995
$ws.Range("M122").Value = -47372.00199999999
$ws.Range("N122").Value = -18371.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4918.1816
$ws.Range("I105").Value = 4350
$ws.Range("J105").Value = 4975
$ws.Range("K105").Value = 4350
$ws.Range("L105").Value = 4975
$ws.Range("M105").Value = -2603
$ws.Range("N105").Value = -8469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5162.047
$ws.Range("I31").Value = 1612.5
$ws.Range("J31").Value = 8294
$ws.Range("K31").Value = 1612.5
$ws.Range("L31").Value = 8294
$ws.Range("M31").Value = -1317.5
$ws.Range("N31").Value = -8884

$ws.Range("H34").Value = 5162.047
$ws.Range("I34").Value = 1612.5
$ws.Range("J34").Value = 8294
$ws.Range("K34").Value = 1612.5
$ws.Range("L34").Value = 8294
$ws.Range("M34").Value = -1410.5
$ws.Range("N34").Value = -8698

$ws.Range("H99").Value = 2633.3333
$ws.Range("I99").Value = 1516.6666
$ws.Range("J99").Value = 3750
$ws.Range("K99").Value = 1516.6666
$ws.Range("L99").Value = 3750
$ws.Range("M99").Value = -18.66660000000002
$ws.Range("N99").Value = -6746

$ws.Range("H126").Value = 2633.3333
$ws.Range("I126").Value = 1516.6666
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 4549.9998
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -2079.9998
$ws.Range("N126").Value = -16190

$ws.Range("H132").Value = 4045.6738
$ws.Range("I132").Value = 1770.6316
$ws.Range("K132").Value = 5311.8948
$ws.Range("M132").Value = -2781.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I12").Value = 44.142857
$ws.Range("J12").Value = 25.2
$ws.Range("K12").Value = 132.428571
$ws.Range("L12").Value = 75.59999999999999
$ws.Range("M12").Value = 40.57142899999999
$ws.Range("N12").Value = -421.6

$ws.Range("H113").Value = 1560.3077
$ws.Range("I113").Value = 1462.75
$ws.Range("K113").Value = 4388.25
$ws.Range("M113").Value = -2218.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18471
$ws.Range("I70").Value = 28974.5
$ws.Range("J70").Value = 4466.3335
$ws.Range("K70").Value = 28974.5
$ws.Range("L70").Value = 4466.3335
$ws.Range("M70").Value = -28704.5
$ws.Range("N70").Value = -5006.3335

$ws.Range("H73").Value = 18471
$ws.Range("I73").Value = 28974.5
$ws.Range("J73").Value = 4466.3335
$ws.Range("K73").Value = 28974.5
$ws.Range("L73").Value = 4466.3335
$ws.Range("M73").Value = -28038.5
$ws.Range("N73").Value = -6338.3335

$ws.Range("H80").Value = 2899
$ws.Range("I80").Value = 2709.4443
$ws.Range("J80").Value = 3088.5557
$ws.Range("K80").Value = 2709.4443
$ws.Range("L80").Value = 3088.5557
$ws.Range("M80").Value = -1711.4443
$ws.Range("N80").Value = -5084.5557

$ws.Range("H83").Value = 2899
$ws.Range("I83").Value = 2709.4443
$ws.Range("J83").Value = 3088.5557
$ws.Range("K83").Value = 13547.2215
$ws.Range("L83").Value = 15442.7785
$ws.Range("M83").Value = -8555.2215
$ws.Range("N83").Value = -25426.7785

$ws.Range("H122").Value = 4128
$ws.Range("I122").Value = 4349
$ws.Range("J122").Value = 4017.5
$ws.Range("K122").Value = 13047
$ws.Range("L122").Value = 12052.5
$ws.Range("M122").Value = -10597
$ws.Range("N122").Value = -16952.5

$ws.Range("H132").Value = 8824.412
$ws.Range("I132").Value = 9063.5
$ws.Range("K132").Value = 27190.5
$ws.Range("M132").Value = -24660.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9691.223
$ws.Range("I122").Value = 18816.666
$ws.Range("J122").Value = 7866.1333
$ws.Range("K122").Value = 56449.99800000001
$ws.Range("L122").Value = 23598.3999
$ws.Range("M122").Value = -53999.99800000001
$ws.Range("N122").Value = -28498.3999
